$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells remain text so numeric-looking values (e.g. "1.00", "6.32")
# are not auto-converted to numbers and keep their exact formatting.

$cells = @("D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "E13", "B14", "C14", "D14", "E14", "B15", "C15", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "E25", "E26", "E27", "D28", "E28", "E29", "D30", "E30", "D31", "E31", "D32", "E32", "E33", "D34", "E34", "B35", "C35", "D35", "E35", "B36", "C36", "D36", "E36", "D37", "E37", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "E44", "D45", "E45", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($cell in $cells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = "39.793.63"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "2.226.38"
$ws.Range("E3").Value = "  -5.14%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "294.41"
$ws.Range("E5").Value = "  -5.20%  "
$ws.Range("D6").Value = "84.50"
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("D7").Value = "0.515"
$ws.Range("E7").Value = "  -2.65%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "0.467"
$ws.Range("E9").Value = "  -3.28%  "
$ws.Range("D10").Value = "0.0788"
$ws.Range("E10").Value = "  -3.07%  "
$ws.Range("D11").Value = "29.92"
$ws.Range("E11").Value = "  -0.70%  "
$ws.Range("D12").Value = "47.83"
$ws.Range("E12").Value = "  -8.82%  "
$ws.Range("E13").Value = "  -2.39%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.572.92"
$ws.Range("E14").Value = "  -4.92%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "6.32"
$ws.Range("E15").Value = "  -1.51%  "
$ws.Range("D16").Value = "14.13"
$ws.Range("E16").Value = "  -4.54%  "
$ws.Range("D17").Value = "2.226.09"
$ws.Range("E17").Value = "  -6.01%  "
$ws.Range("D18").Value = "0.721"
$ws.Range("E18").Value = "  -5.27%  "
$ws.Range("D19").Value = "39.731.24"
$ws.Range("E19").Value = "  -0.85%  "
$ws.Range("D20").Value = "0.0₃0884"
$ws.Range("E20").Value = "  -2.01%  "
$ws.Range("D21").Value = "5.77"
$ws.Range("E21").Value = "  -5.49%  "
$ws.Range("D22").Value = "65.33"
$ws.Range("E22").Value = "  -4.24%  "
$ws.Range("D23").Value = "10.51"
$ws.Range("E23").Value = "  -1.40%  "
$ws.Range("D24").Value = "232.46"
$ws.Range("E24").Value = "  -1.19%  "
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("E26").Value = "  -5.33%  "
$ws.Range("E27").Value = "  +0.70%  "
$ws.Range("D28").Value = "22.86"
$ws.Range("E28").Value = "  -4.01%  "
$ws.Range("E29").Value = "  -0.75%  "
$ws.Range("D30").Value = "9.20"
$ws.Range("E30").Value = "  -0.97%  "
$ws.Range("D31").Value = "32.45"
$ws.Range("E31").Value = "  -6.48%  "
$ws.Range("D32").Value = "151.08"
$ws.Range("E32").Value = "  -1.76%  "
$ws.Range("E33").Value = "  -0.14%  "
$ws.Range("D34").Value = "4.82"
$ws.Range("E34").Value = "  -5.77%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "0.0704"
$ws.Range("E35").Value = "  -1.97%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "2.37"
$ws.Range("E36").Value = "  -3.83%  "
$ws.Range("D37").Value = "16.04"
$ws.Range("E37").Value = "  +2.86%  "
$ws.Range("E38").Value = "  -1.91%  "
$ws.Range("D39").Value = "0.0979"
$ws.Range("E39").Value = "  -0.55%  "
$ws.Range("D40").Value = "2.66"
$ws.Range("E40").Value = "  -5.20%  "
$ws.Range("D41").Value = "1.65"
$ws.Range("E41").Value = "  -4.08%  "
$ws.Range("D42").Value = "3.70"
$ws.Range("E42").Value = "  -3.86%  "
$ws.Range("D43").Value = "1.947.86"
$ws.Range("E43").Value = "  -0.88%  "
$ws.Range("E44").Value = "  -3.75%  "
$ws.Range("D45").Value = "0.0267"
$ws.Range("E45").Value = "  +0.88%  "
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("D47").Value = "16.21"
$ws.Range("E47").Value = "  -7.85%  "
$ws.Range("D48").Value = "2.59"
$ws.Range("E48").Value = "  -4.17%  "
$ws.Range("D49").Value = "2.440.25"
$ws.Range("E49").Value = "  -4.92%  "
$ws.Range("D50").Value = "70.68"
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("D51").Value = "88.99"
$ws.Range("E51").Value = "  -4.54%  "
